$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The column C ("caseDetailQuery" header + its Cypher query text in row 2)
# is no longer needed. Delete the entire column C, which shifts the
# former columns D (dbExcel / Neo4jData filename) and E (WebExcel / WebData
# filename) left to become the new C and D.
$ws.Range("C:C").Delete()

# Keep the selection consistent with what Excel leaves behind after a
# whole-column delete (the entire new column C is selected).
$ws.Range("C:C").Select()
